$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 616, shifting existing rows 616-659 down to 617-660
$ws.Rows(616).Insert()

# Populate the newly inserted row 616 with the new data record
$ws.Range("A616").Value = 9
$ws.Range("B616").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C616").Value = "Metropolitana"
$ws.Range("D616").Value = 45265
$ws.Range("E616").Value = 13
$ws.Range("F616").Value = 100112052
$ws.Range("G616").Value = "Albahaca"
$ws.Range("H616").Value = "Sin especificar"
$ws.Range("I616").Value = "Primera"
$ws.Range("J616").Value = 250
$ws.Range("K616").Value = 6000
$ws.Range("L616").Value = 6000
$ws.Range("M616").Value = 6000
$ws.Range("N616").Value = "$/docena de matas"
$ws.Range("O616").Value = "Provincia de Chacabuco"
$ws.Range("P616").Value = 1000
$ws.Range("Q616").Value = 6
$ws.Range("R616").Value = "Hortaliza"
